# Apply the two changes captured by the commit:
#
#  1. Slide 6's table switches from the embedded custom "Table_0" style
#     ({EE95E160-03D7-4FDE-9BAA-573EE347DD88}) to the built-in table
#     style {1A027DD0-B02C-4D27-B786-AE598EF23C9C}.
#
#  2. The deck's theme colour palette is swapped from the "Integral"
#     palette back to the stock "Office" palette (dk1/lt1/dk2/lt2/
#     accent1-6/hlink/folHlink), matching the colours that the
#     presentation's secondary ("Office Theme") theme part used before
#     the edit.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s6 = $p.Slides.Item(6)
$tbl = $s6.Shapes.Item(2).Table
$tbl.ApplyStyle("{1A027DD0-B02C-4D27-B786-AE598EF23C9C}")

# --- 2. Theme colours -------------------------------------------------
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $rr = ($hex -shr 16) -band 0xFF
    $gg = ($hex -shr 8) -band 0xFF
    $bb = $hex -band 0xFF
    $bgr = $rr + ($gg * 256) + ($bb * 65536)
    $tcs.Item($i).RGB = $bgr
}
